# DPLKINV001-011 - Setup Jenis Emiten Investasi
# Replace the long, numbered step-by-step instruction text in column D
# (SCENARIO_DESC) with short "<Action> Setup Jenis Emiten" labels, adjust
# the affected row heights, and move the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Tambah :..." -> "Tambah Setup Jenis Emiten"
$ws.Range("D2").Value = "Tambah Setup Jenis Emiten"

# Row 3: "View :..." -> "View Setup Jenis Emiten"
$ws.Range("D3").Value = "View Setup Jenis Emiten"

# Row 4: "Ubah :..." -> "Ubah Setup Jenis Emiten"
$ws.Range("D4").Value = "Ubah Setup Jenis Emiten"

# Row 5: "Hapus :..." -> "Hapus Setup Jenis Emiten"
$ws.Range("D5").Value = "Hapus Setup Jenis Emiten"

# The shortened text needs fewer wrapped lines, so shrink the row heights
# accordingly (rows 2, 3 and 5 end up at 30pt; row 4 returns to the sheet's
# default/auto height).
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 30

# Move the active selection to D5.
$ws.Range("D5").Select()
